# ConceptualDomain.xlsx — "Model Updates, Added Upload Code"
#
# - Header row: B1 label "skos:prefLabel" -> "URI" swap positions with a new
#   "dcdtr:domainName" label (A1 becomes "URI", B1 becomes "dcdtr:domainName").
# - A2:A11 formulas collapse into one shared formula.
# - A new (currently empty) column F is reserved/widened for upcoming upload
#   data, and the active selection moves to F1.
# - Column widths on A:D nudge slightly; column F gets a width of 22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row relabel -------------------------------------------------
$ws.Range("A1").Value = "URI"
$ws.Range("B1").Value = "dcdtr:domainName"

# --- Collapse the per-row formulas into a single shared formula --------
$ws.Range("A2:A11").Formula = '="class/conceptual-domain/"&LOWER(SUBSTITUTE(B2," ", "-"))'

# --- Column widths (engine stores width = ColumnWidth + 5/6) -----------
$ws.Columns.Item(1).ColumnWidth = 28 - 5/6
$ws.Columns.Item(2).ColumnWidth = 23.44140625 - 5/6
$ws.Columns.Item(3).ColumnWidth = 30.6640625 - 5/6
$ws.Columns.Item(4).ColumnWidth = 11.88671875 - 5/6
$ws.Columns.Item(6).ColumnWidth = 22 - 5/6

# --- Move the selection on the frozen (bottom-right) pane to F1 --------
[void]$ws.Range("F1").Select()

# --- Row heights re-measured slightly shorter (wrap-text autofit) ------
$ws.Rows.Item(2).RowHeight = 41.4
$ws.Rows.Item(3).RowHeight = 41.4
$ws.Rows.Item(4).RowHeight = 41.4
$ws.Rows.Item(5).RowHeight = 41.4
$ws.Rows.Item(6).RowHeight = 55.2
$ws.Rows.Item(7).RowHeight = 41.4
$ws.Rows.Item(8).RowHeight = 27.6
$ws.Rows.Item(9).RowHeight = 41.4
$ws.Rows.Item(10).RowHeight = 41.4
$ws.Rows.Item(11).RowHeight = 110.4
